$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.979.44'
$ws.Range('E2').Value = '  -2.50%  '
$ws.Range('D3').Value = '3.186.92'
$ws.Range('E3').Value = '  -1.44%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.37'
$ws.Range('E5').Value = '  -1.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.12'
$ws.Range('E6').Value = '  -3.36%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '3.185.64'
$ws.Range('E8').Value = '  -1.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.534'
$ws.Range('E9').Value = '  -3.45%  '
$ws.Range('E10').Value = '  -4.16%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.61'
$ws.Range('E11').Value = '  -1.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.480'
$ws.Range('E12').Value = '  -5.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000263'
$ws.Range('E13').Value = '  -3.65%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.32'
$ws.Range('E14').Value = '  -4.17%  '
$ws.Range('D15').Value = '3.704.44'
$ws.Range('E15').Value = '  -1.59%  '
$ws.Range('D16').Value = '64.956.83'
$ws.Range('E16').Value = '  -2.63%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.114'
$ws.Range('E17').Value = '  +0.68%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.168.53'
$ws.Range('E18').Value = '  -1.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.11'
$ws.Range('E19').Value = '  -3.91%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '486.67'
$ws.Range('E20').Value = '  -4.82%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.92'
$ws.Range('E21').Value = '  -1.89%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.722'
$ws.Range('E22').Value = '  -2.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.85'
$ws.Range('E23').Value = '  -2.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.10'
$ws.Range('E24').Value = '  -3.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.13'
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.95'
$ws.Range('E27').Value = '  -2.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.83'
$ws.Range('E28').Value = '  -2.99%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.28'
$ws.Range('E29').Value = '  -4.54%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.22'
$ws.Range('E30').Value = '  +2.61%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.122'
$ws.Range('E31').Value = '  +5.91%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.74'
$ws.Range('E32').Value = '  -7.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.08'
$ws.Range('E33').Value = '  -3.99%  '
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('E35').Value = '  -5.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.19'
$ws.Range('E36').Value = '  -4.81%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '54.64'
$ws.Range('E37').Value = '  -1.47%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.25'
$ws.Range('E38').Value = '  +6.54%  '
$ws.Range('D39').Value = '0.0₃0744'
$ws.Range('E39').Value = '  -3.75%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '465.70'
$ws.Range('E40').Value = '  -8.17%  '
$ws.Range('E41').Value = '  -1.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0408'
$ws.Range('E42').Value = '  -3.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.59'
$ws.Range('E43').Value = '  -1.87%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.48'
$ws.Range('E44').Value = '  +1.60%  '
$ws.Range('D45').Value = '2.935.62'
$ws.Range('E45').Value = '  +1.19%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.279'
$ws.Range('E46').Value = '  -6.67%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '27.56'
$ws.Range('E47').Value = '  -2.46%  '
$ws.Range('B48').Value = 'ThetaToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.38'
$ws.Range('E48').Value = '  -1.95%  '
$ws.Range('B49').Value = 'USDe'
$ws.Range('C49').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.999'
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.93'
$ws.Range('E51').Value = '  -2.17%  '
